$d = $word.ActiveDocument

# Locate the paragraph that starts the "Mata Garuda merupakan ..." block.
$p31 = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text.StartsWith("Mata Garuda merupakan ")) {
        $p31 = $i
        break
    }
}

$r = $d.Paragraphs($p31).Range
$found = $r.Find.Execute("yang memonitor*harinya.", $true, $false, $true, $false, $false, $true, 1, $false, "", 0)
$newText = "yang memonitor lalu lintas data pada jaringan internet yang ada di Indonesia. Mata Garuda akan melaporkan suatu kejadian yang cirinya telah didefinisikan dalam sebuah rule. Kejadian-kejadian tersebut didapatkan melalui sensor yang terpasang di netiap Network Access Point yang ada di Indonesia.  Sensor tersebut berfungsi untuk mengambil paket lalu meneruskannya ke defense center.  Setiap harinya terdapat puluhan juta packet yang ditangkap oleh sensor Mata Garuda dan dengan terdeteksi rata-rata 2 juta serangan dideteksi per harinya."
$r.Text = $newText

$pr = $d.Paragraphs($p31).Range
$null = $pr.Find.Execute("Network Access Point", $false, $true, $false, $false, $false, $true, 1, $false, "", 0)
$pr.Italic = 1

$pr2 = $d.Paragraphs($p31).Range
$null = $pr2.Find.Execute("defense center", $false, $true, $false, $false, $false, $true, 1, $false, "", 0)
$pr2.Italic = 1

Write-Host $d.Paragraphs($p31).Range.Text

# --- Paragraph 2: "Dengan Semakin berkembangnya..." ---
$p32 = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text.StartsWith("Dengan Semakin berkembangnya")) {
        $p32 = $i
        break
    }
}

$r2 = $d.Paragraphs($p32).Range
$null = $r2.Find.Execute(". Dengan semakin*tersebut.", $true, $false, $true, $false, $false, $true, 1, $false, "", 0)
$r2.Text = ". Dengan semakin besarnya lalu lintas data pada jaringan akan menyebabkan semakin banyaknya event yang harus diproses dan dianalisa oleh defense center Mata Garuda. Hal ini berpengaruh secara langsung dalam kecepatan Mata Garuda dalam menganlisa data lalu lintas tersebut. Oleh karena itu, pengembangan arsitektur sistem Mata Garuda harus dilakukan agar dapat sesuai dengan kondisi sekarang"

$pr3 = $d.Paragraphs($p32).Range
$null = $pr3.Find.Execute("event", $false, $true, $false, $false, $false, $true, 1, $false, "", 0)
$pr3.Italic = 1

$pr4 = $d.Paragraphs($p32).Range
$null = $pr4.Find.Execute("defense center", $false, $true, $false, $false, $false, $true, 1, $false, "", 0)
$pr4.Italic = 1

Write-Host $d.Paragraphs($p32).Range.Text

# --- Remove the old third paragraph ("Untuk menganalisa data yang sangat besar...") ---
$p33 = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text.StartsWith("Untuk menganalisa data yang sangat besar")) {
        $p33 = $i
        break
    }
}
$d.Paragraphs($p33).Range.Delete()

# --- Remove the now-empty paragraph that followed it ---
$pEmpty = $d.Paragraphs($p33)
if ($pEmpty.Range.Text.Trim().Length -eq 0) {
    $pEmpty.Range.Delete()
}

# --- Bookmark "_GoBack" spanning from the start of paragraph 1 to the end
#     (before the paragraph mark) of paragraph 2 ---
$startPos = $d.Paragraphs($p31).Range.Start
$endRange = $d.Paragraphs($p32).Range
$endPos = $endRange.End - 1
$bmRange = $d.Range($startPos, $endPos)
Write-Host "bookmark text = [$($bmRange.Text)]"
$d.Bookmarks.Add("_GoBack", $bmRange)
